# Refresh computed market-price / profit columns (H:N) across several
# worksheets with updated Universalis price data.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 15
$ws.Range("H15").Value = 1432415.8
$ws.Range("I15").Value = 1432415.8
$ws.Range("K15").Value = 4297247.4
$ws.Range("M15").Value = -4297078.4

# row 33
$ws.Range("H33").Value = 545.4138
$ws.Range("I33").Value = 682.6818
$ws.Range("J33").Value = 114
$ws.Range("K33").Value = 682.6818
$ws.Range("L33").Value = 114
$ws.Range("M33").Value = -453.6818
$ws.Range("N33").Value = -572

# row 80
$ws.Range("H80").Value = 568.2143
$ws.Range("I80").Value = 703.3333
$ws.Range("J80").Value = 504.21054
$ws.Range("K80").Value = 2109.9999
$ws.Range("L80").Value = 1512.63162
$ws.Range("M80").Value = -1111.9999
$ws.Range("N80").Value = -3508.63162

# row 83
$ws.Range("H83").Value = 568.2143
$ws.Range("I83").Value = 703.3333
$ws.Range("J83").Value = 504.21054
$ws.Range("K83").Value = 6329.9997
$ws.Range("L83").Value = 4537.894859999999
$ws.Range("M83").Value = -1337.9997
$ws.Range("N83").Value = -14521.89486

# row 135
$ws.Range("H135").Value = 742.9318
$ws.Range("I135").Value = 280.38095
$ws.Range("K135").Value = 2523.42855
$ws.Range("M135").Value = 11.57145000000037

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 55
$ws.Range("H55").Value = 29000
$ws.Range("J55").Value = 29000
$ws.Range("L55").Value = 29000
$ws.Range("N55").Value = -29546

# row 135
$ws.Range("H135").Value = 37780
$ws.Range("J135").Value = 37780
$ws.Range("L135").Value = 37780
$ws.Range("N135").Value = -47920

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 7
$ws.Range("H7").Value = 5938.5884
$ws.Range("I7").Value = 16711.834
$ws.Range("J7").Value = 62.272728
$ws.Range("K7").Value = 16711.834
$ws.Range("L7").Value = 62.272728
$ws.Range("M7").Value = -16598.834
$ws.Range("N7").Value = -288.272728

# row 22
$ws.Range("H22").Value = 340.57895
$ws.Range("I22").Value = 248.41667
$ws.Range("J22").Value = 498.57144
$ws.Range("K22").Value = 248.41667
$ws.Range("L22").Value = 498.57144
$ws.Range("M22").Value = 101.58333
$ws.Range("N22").Value = -1198.57144

# row 31
$ws.Range("H31").Value = 5466246.5
$ws.Range("I31").Value = 1431.7407
$ws.Range("J31").Value = 47623388
$ws.Range("K31").Value = 1431.7407
$ws.Range("L31").Value = 47623388
$ws.Range("M31").Value = -1136.7407
$ws.Range("N31").Value = -47623978

# row 34
$ws.Range("H34").Value = 5466246.5
$ws.Range("I34").Value = 1431.7407
$ws.Range("J34").Value = 47623388
$ws.Range("K34").Value = 1431.7407
$ws.Range("L34").Value = 47623388
$ws.Range("M34").Value = -1229.7407
$ws.Range("N34").Value = -47623792

# row 58
$ws.Range("H58").Value = 4370.425
$ws.Range("I58").Value = 3038.389
$ws.Range("J58").Value = 5460.273
$ws.Range("K58").Value = 3038.389
$ws.Range("L58").Value = 5460.273
$ws.Range("M58").Value = -2835.389
$ws.Range("N58").Value = -5866.273

# row 122
$ws.Range("H122").Value = 58825490
$ws.Range("I122").Value = 90910300
$ws.Range("J122").Value = 3324.1667
$ws.Range("K122").Value = 272730900
$ws.Range("L122").Value = 9972.500100000001
$ws.Range("M122").Value = -272728450
$ws.Range("N122").Value = -14872.5001

# row 136
$ws.Range("H136").Value = 4370.425
$ws.Range("I136").Value = 3038.389
$ws.Range("J136").Value = 5460.273
$ws.Range("K136").Value = 9115.167000000001
$ws.Range("L136").Value = 16380.819
$ws.Range("M136").Value = -6565.167000000001
$ws.Range("N136").Value = -21480.819

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 5
$ws.Range("H5").Value = 542.95654
$ws.Range("I5").Value = 253.64706
$ws.Range("J5").Value = 1362.6666
$ws.Range("K5").Value = 760.94118
$ws.Range("L5").Value = 4087.9998
$ws.Range("M5").Value = -648.94118
$ws.Range("N5").Value = -4311.9998

# row 107
$ws.Range("H107").Value = 27778164
$ws.Range("I107").Value = 71428940
$ws.Range("J107").Value = 394.9091
$ws.Range("K107").Value = 214286820
$ws.Range("L107").Value = 1184.7273
$ws.Range("M107").Value = -214284900
$ws.Range("N107").Value = -5024.7273

# row 113
$ws.Range("H113").Value = 522.2406999999999
$ws.Range("I113").Value = 391.66666
$ws.Range("J113").Value = 783.3889
$ws.Range("K113").Value = 1174.99998
$ws.Range("L113").Value = 2350.1667
$ws.Range("M113").Value = 995.0000199999999
$ws.Range("N113").Value = -6690.1667

# row 132
$ws.Range("H132").Value = 974.3333
$ws.Range("I132").Value = 895.86664
$ws.Range("J132").Value = 1366.6666
$ws.Range("K132").Value = 8062.79976
$ws.Range("L132").Value = 12299.9994
$ws.Range("M132").Value = -5532.79976
$ws.Range("N132").Value = -17359.9994

# row 135
$ws.Range("H135").Value = 542.95654
$ws.Range("I135").Value = 253.64706
$ws.Range("J135").Value = 1362.6666
$ws.Range("K135").Value = 2282.82354
$ws.Range("L135").Value = 12263.9994
$ws.Range("M135").Value = 252.1764599999997
$ws.Range("N135").Value = -17333.9994

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 1700
$ws.Range("I22").Value = 1250
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 1250
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -955
$ws.Range("N22").Value = -2590

# row 27
$ws.Range("H27").Value = 1700
$ws.Range("I27").Value = 1250
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 1250
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -1143
$ws.Range("N27").Value = -2214

# row 40
$ws.Range("H40").Value = 27029068
$ws.Range("I40").Value = 37038684
$ws.Range("J40").Value = 3108
$ws.Range("K40").Value = 37038684
$ws.Range("L40").Value = 3108
$ws.Range("M40").Value = -37038548
$ws.Range("N40").Value = -3380

# row 93
$ws.Range("H93").Value = 1945.6
$ws.Range("I93").Value = 1900.5625
$ws.Range("J93").Value = 2125.75
$ws.Range("K93").Value = 1900.5625
$ws.Range("L93").Value = 2125.75
$ws.Range("M93").Value = -652.5625
$ws.Range("N93").Value = -4621.75

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 126
$ws.Range("H126").Value = 1556.4902
$ws.Range("I126").Value = 1810.1578
$ws.Range("J126").Value = 815
$ws.Range("K126").Value = 5430.4734
$ws.Range("L126").Value = 2445
$ws.Range("M126").Value = -2960.4734
$ws.Range("N126").Value = -7385

# row 129
$ws.Range("H129").Value = 29800
$ws.Range("J129").Value = 29800
$ws.Range("L129").Value = 29800
$ws.Range("N129").Value = -39800

# row 130
$ws.Range("H130").Value = 29444.375
$ws.Range("J130").Value = 29444.375
$ws.Range("L130").Value = 29444.375
$ws.Range("N130").Value = -39484.375
